$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-09-30 00:00:00"
$ws.Range("O2").Value = 202869167.13
$ws.Range("P2").Value = 1364786417.99
$ws.Range("Q2").Value = 1153579081.36
$ws.Range("R2").Value = 25.6355642712
$ws.Range("S2").Value = 962027431.8099999
$ws.Range("T2").Value = 962027431.8099999
$ws.Range("U2").Value = 38.6848679984
$ws.Range("V2").Value = 47729124.05
$ws.Range("W2").Value = 84887333.09999999
$ws.Range("X2").Value = 15472359.61
$ws.Range("Y2").Value = 274187759.31
$ws.Range("Z2").Value = 274255101.13
$ws.Range("AA2").Value = 71385934
$ws.Range("AG2").Value = 43462832.79
$ws.Range("AP2").Value = 28.2287060491
$ws.Range("AQ2").Value = -82.803618622306
$ws.Range("AR2").Value = -83.79788318052501
$ws.Range("AS2").Value = 152633654.2
$ws.Range("AT2").Value = 55.63777918967
